$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily price-data refresh: insert a new "latest" row just under the header
# (row 2), pushing every existing date row down by one. The new row carries
# the same price figures as the prior newest row (783.5 / 1112 / 3610) but
# dated one day later.
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-01-01"
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610

# Inserting a row above copies formatting down from the header row (bold,
# centered, bordered); the new data row should look like the other plain
# data rows instead, so drop that inherited formatting.
$ws.Rows.Item(2).ClearFormats()
